$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create "basebert1" as a duplicate of "basebert" (right after it) so it
# inherits the exact same sheet properties / formatting / styles, then
# overwrite its data with the new classifier run's results.
# ---------------------------------------------------------------------------
$baseSheet = $wb.Worksheets.Item("basebert")
$baseSheet.Copy($null, $baseSheet)
$newSheet = $wb.Worksheets.Item("basebert (2)")
$newSheet.Name = "basebert1"

# ---------------------------------------------------------------------------
# Data rows - index(A), fold(B), epoch(D) and the numeric metrics (I-N) are
# plain values; version(C) and recall/precision/accuracy/fbeta (E-H) are
# stored as TEXT, same as in the sibling "basebert" sheet this was copied
# from, so plain .Value assignment keeps them text automatically.
# ---------------------------------------------------------------------------
$data = @(
    @(0, "fold_0", "18.09_09.50", 8, "0.71794873", "0.23140496", "0.82038",   "0.50541514", 0, 0, 0.5, 22, 186, 0.5986553269128004),
    @(1, "fold_1", "18.09_10.00", 6, "0.8076923",  "0.23684211", "0.8117444", "0.5449827",  0, 0, 0.5, 15, 203, 0.7180119156837463),
    @(2, "fold_2", "18.09_10.12", 5, "0.8717949",  "0.32380953", "0.8687392", "0.651341",   0, 0, 0.5, 10, 142, 0.8751945644617081),
    @(3, "fold_3", "18.09_10.25", 8, "0.8101266",  "0.24521072", "0.8169257", "0.5545927",  0, 0, 0.5, 15, 197, 0.8309154734015465)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $rowIndex = $r + 2

    $newSheet.Cells.Item($rowIndex, 1).Value  = $row[0]   # A - index
    $newSheet.Cells.Item($rowIndex, 2).Value  = $row[1]   # B - fold
    $newSheet.Cells.Item($rowIndex, 3).Value  = $row[2]   # C - version
    $newSheet.Cells.Item($rowIndex, 4).Value  = $row[3]   # D - epoch

    # E:H hold numeric-looking text ("0.71794873", ...) - force a Text
    # number format before assigning so the engine keeps them as strings
    # instead of silently parsing them into doubles.
    $ehRange = $newSheet.Range("E${rowIndex}:H${rowIndex}")
    $ehRange.NumberFormat = "@"
    $newSheet.Cells.Item($rowIndex, 5).Value  = $row[4]   # E - recall
    $newSheet.Cells.Item($rowIndex, 6).Value  = $row[5]   # F - precision
    $newSheet.Cells.Item($rowIndex, 7).Value  = $row[6]   # G - accuracy
    $newSheet.Cells.Item($rowIndex, 8).Value  = $row[7]   # H - fbeta

    $newSheet.Cells.Item($rowIndex, 9).Value  = $row[8]   # I - best recall
    $newSheet.Cells.Item($rowIndex, 10).Value = $row[9]   # J - best precision
    $newSheet.Cells.Item($rowIndex, 11).Value = $row[10]  # K - best threshold
    $newSheet.Cells.Item($rowIndex, 12).Value = $row[11]  # L - false neg
    $newSheet.Cells.Item($rowIndex, 13).Value = $row[12]  # M - false pos
    $newSheet.Cells.Item($rowIndex, 14).Value = $row[13]  # N - val loss
}

# Restore the plain (unstyled / General-format) look on E2:H5 now that the
# values are stored as text - copy the formatting from the untouched D
# column (still default/General) across so no visible number-format change
# remains on the new sheet.
$newSheet.Range("D2:D5").Copy()
$newSheet.Range("E2").PasteSpecial(-4122)
$newSheet.Range("F2").PasteSpecial(-4122)
$newSheet.Range("G2").PasteSpecial(-4122)
$newSheet.Range("H2").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
